$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Metadata sheet: bump the "Date" value (row 8, column B).
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2024-03-22T16:25:12+00:00"

# ---------------------------------------------------------------------------
# 2) Elements sheet: swap the two "Mapping" columns (AK <-> AL), header and
#    data alike, for rows 1-14 (the used range of the sheet).
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Elements")

for ($r = 1; $r -le 14; $r++) {
    $akCell = $ws.Cells.Item($r, 37)
    $alCell = $ws.Cells.Item($r, 38)

    $akText = $akCell.Text
    $alText = $alCell.Text

    # Only touch cells whose value actually needs to change - this avoids
    # needlessly clobbering cells that are identical before/after the swap
    # (e.g. rows where both columns are already blank).
    if ($akText -ne $alText) {
        $akCell.Value = $alText
        $alCell.Value = $akText
    }
}

# Swap the column widths that went along with the two swapped columns.
$ws.Columns.Item(37).ColumnWidth = 89.9296875
$ws.Columns.Item(38).ColumnWidth = 24.98046875
